$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.292.57"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.42%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.490.47"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.76"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.75%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.15"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.04%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.72"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +4.83%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.31%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.123"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.17%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.32"
$ws.Range("D13").Style = $style
$ws.Range("E14").Value = "  -1.34%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.880.62"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.35%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.491.18"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -0.95%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.36%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.190.20"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.37%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.92"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +2.44%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.60"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.07%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0932"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +12.84%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.25"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  -2.02%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +0.00%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.73"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.55"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  +0.38%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.76"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +5.92%  "
$ws.Range("E34").Value = "  +0.46%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.72"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +1.70%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.92"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -1.18%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.11"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +8.61%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  +0.09%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.05"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("E44").Value = "  +0.10%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.989.41"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.76%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  -5.70%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  -5.67%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.62"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +4.17%  "
